# Update "想去人数" (want-to-go count) values that changed between data refreshes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 7347
$wsExhibit.Range("F6").Value = 435
$wsExhibit.Range("F7").Value = 3867

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 7347
$wsAll.Range("F8").Value = 435
$wsAll.Range("F9").Value = 3867
